$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2896
$ws1.Range("F3").Value = 21255
$ws1.Range("F4").Value = 104
$ws1.Range("F5").Value = 3095
$ws1.Range("F6").Value = 812
$ws1.Range("F7").Value = 620
$ws1.Range("F9").Value = 774
$ws1.Range("F10").Value = 285
$ws1.Range("F12").Value = 73
$ws1.Range("F14").Value = 529
$ws1.Range("F15").Value = 185
$ws1.Range("F16").Value = 276
$ws1.Range("F17").Value = 21
$ws1.Range("F18").Value = 426
$ws1.Range("F19").Value = 72
$ws1.Range("F22").Value = 42

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F20").Value = 39

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6137
$ws3.Range("F4").Value = 702
$ws3.Range("F5").Value = 1647
$ws3.Range("F6").Value = 58

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6137
$ws4.Range("F4").Value = 702
$ws4.Range("F5").Value = 1647
$ws4.Range("F6").Value = 2896
$ws4.Range("F7").Value = 21255
$ws4.Range("F10").Value = 104
$ws4.Range("F13").Value = 3095
$ws4.Range("F14").Value = 812
$ws4.Range("F16").Value = 58
$ws4.Range("F17").Value = 620
$ws4.Range("F19").Value = 774
$ws4.Range("F20").Value = 285
$ws4.Range("F23").Value = 73
$ws4.Range("F29").Value = 529
$ws4.Range("F31").Value = 185
$ws4.Range("F33").Value = 276
$ws4.Range("F36").Value = 21
$ws4.Range("F37").Value = 426
$ws4.Range("F39").Value = 72
$ws4.Range("F44").Value = 42
$ws4.Range("F47").Value = 39
